# "Add files via upload" — append the next day's data row to the two
# pool sheets ("gUSD 26.06.25" and "mPendle 27.03.25") and update the
# selection / active-sheet UI state to match where the author ended up.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# gUSD 26.06.25  (sheet "gUSD 26.06.25") — new data for row 7
# ---------------------------------------------------------------------
$wsGUSD0626 = $wb.Worksheets.Item("gUSD 26.06.25")

$wsGUSD0626.Range("B7").Value = 0.036889999999999999
$wsGUSD0626.Range("C7").Value = 42.27
$wsGUSD0626.Range("D7").Value = 0.76619999999999999
$wsGUSD0626.Range("E7").Value = 7.93
$wsGUSD0626.Range("F7").Value = 8.0399999999999991
$wsGUSD0626.Range("G7").Value = 9.59

# ---------------------------------------------------------------------
# mPendle 27.03.25 (sheet "mPendle 27.03.25") — new data for row 6
# ---------------------------------------------------------------------
$wsPendle = $wb.Worksheets.Item("mPendle 27.03.25")

# The B column on this sheet is formatted with one extra decimal place
# (0.00000) compared to the default (0.0000) used elsewhere.
$wsPendle.Range("B5:B6").NumberFormat = "0.00000"

$wsPendle.Range("B6").Value = 0.016670000000000001
$wsPendle.Range("C6").Value = 25.68
$wsPendle.Range("D6").Value = 0.43109999999999998
$wsPendle.Range("E6").Value = 4.74
$wsPendle.Range("F6").Value = 5.38
$wsPendle.Range("G6").Value = 4.51

# ---------------------------------------------------------------------
# Selection / active sheet bookkeeping, matching where the author left
# the cursor after entering the new rows.
# ---------------------------------------------------------------------
$wsGUSD0626.Range("G8").Select() | Out-Null

$wsPendle.Range("N6").Select() | Out-Null
$wsPendle.Activate()
